$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.540.20'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '2.023.95'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''254.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.35%  '
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''56.68'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -8.40%  '
$ws.Range("D9").Value = '''0.381'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("D11").Value = '''0.102'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.21%  '
$ws.Range("D12").Value = '''14.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.52%  '
$ws.Range("D13").Value = '2.322.26'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("D14").Value = '''0.813'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.02%  '
$ws.Range("D15").Value = '''21.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.47%  '
$ws.Range("E16").Value = '  -2.60%  '
$ws.Range("D17").Value = '2.026.67'
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '37.422.01'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("E20").Value = '  -2.62%  '
$ws.Range("D21").Value = '''5.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("D22").Value = '''228.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").Value = '''2.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.46%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '''2.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.57%  '
$ws.Range("D26").Value = '''163.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("D27").Value = '''9.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.86%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.132'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.79%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '''19.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("D32").Value = '''0.0668'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.39%  '
$ws.Range("D33").Value = '''4.68'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.38%  '
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("D35").Value = '''2.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.67%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("D39").Value = '''5.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.80%  '
$ws.Range("E40").Value = '  +2.93%  '
$ws.Range("D41").Value = '''0.0962'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.26%  '
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").Value = '1.410.97'
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''90.54'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''15.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.90%  '
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("D50").Value = '''2.02'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").Value = '2.214.19'
$ws.Range("E51").Value = '  +0.99%  '
